$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 2.3
$ws.Range("L2").Value = 6.5
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("AE2").Value = 23
$ws.Range("AH2").Value = 29
$ws.Range("AQ2").Value = 29
$ws.Range("AU2").Value = 10
$ws.Range("AZ2").Value = 151
$ws.Range("BB2").Value = 501

# Row 4
$ws.Range("G4").Value = 3.1
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 2.38
$ws.Range("J4").Value = 3.75
$ws.Range("N4").Value = 9
$ws.Range("AA4").Value = 26
$ws.Range("AH4").Value = 11
$ws.Range("AK4").Value = 21
$ws.Range("AL4").Value = 34
$ws.Range("AQ4").Value = 51

# Row 7
$ws.Range("G7").Value = 2.38
$ws.Range("I7").Value = 3.2
$ws.Range("J7").Value = 3.2
$ws.Range("L7").Value = 4
$ws.Range("O7").Value = 1.5
$ws.Range("P7").Value = 2.5

# Row 10
$ws.Range("G10").Value = 3.3
$ws.Range("J10").Value = 4
$ws.Range("N10").Value = 9.5
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 3.25
$ws.Range("R10").Value = 1.75
$ws.Range("S10").Value = 1.44
$ws.Range("T10").Value = 2.63
$ws.Range("U10").Value = 1.83
$ws.Range("V10").Value = 1.83
$ws.Range("AA10").Value = 29
$ws.Range("AC10").Value = 9.5
$ws.Range("AN10").Value = 5.5
$ws.Range("AP10").Value = 29
$ws.Range("AQ10").Value = 67
$ws.Range("AT10").Value = 2.63
$ws.Range("AW10").Value = 4
$ws.Range("AY10").Value = 23
